$wb = $excel.ActiveWorkbook

# --- Rename "Cancelar encomenda" -> "UC 1 - Registar Organização" ---
$ws = $wb.Worksheets.Item("Cancelar encomenda")
$ws.Name = "UC 1 - Registar Organização"

# --- Row 9: Test Scenario Description ---
$scenarioDesc = "Check Registar Organização functionality"
$ws.Range("B9").Value = $scenarioDesc
$ws.Range("C9").Value = $scenarioDesc
$ws.Range("D9").Value = $scenarioDesc
$ws.Range("E9").Value = $scenarioDesc
$ws.Rows.Item(9).RowHeight = 31.5

# --- Row 10: Test Case ID ---
$ws.Range("B10").Value = "TC_UC1_001"
$ws.Range("C10").Value = "TC_UC1_002"
$ws.Range("D10").Value = "TC_UC1_003"
$ws.Range("E10").Value = "TC_UC1_004"

# --- Row 11: Test Case Description ---
$tcDesc = "Valid Org Name`nValid Org NIF`nValid Org Email`nValid Org Website`nValid Org Phone Number`nValid Org Address"
$ws.Range("B11").Value = $tcDesc
$ws.Range("C11").Value = $tcDesc
$ws.Range("D11").Value = $tcDesc
$ws.Range("E11").Value = $tcDesc
$ws.Rows.Item(11).RowHeight = 94.5

# --- Row 12: Test Case Steps ---
$tcSteps = "1. Enter valid org name`n2. Enter valid org nif`n3. Enter valid org email`n4. Enter valid org website`n5. Enter valid org phone number`n6. Enter valid org address`n7. Click register button"
$ws.Range("B12").Value = $tcSteps
$ws.Range("C12").Value = $tcSteps
$ws.Range("D12").Value = $tcSteps
$ws.Range("E12").Value = $tcSteps
$ws.Rows.Item(12).RowHeight = 126

# --- Row 13: Preconditions (matches source text literally) ---
$ws.Range("B13").Value = "Test Data"
$ws.Range("C13").Value = "Test Data"
$ws.Range("D13").Value = "Test Data"
$ws.Range("E13").Value = "Test Data"

# --- Row 14: Test Data ---
$testData = "org name: org123`norg nif: 123456789`norg email: org123@org123.com`norg website: org123.com`norg phone number: 911123123`norg address: rua sousa, 1, 4000-400, Porto"
$ws.Range("B14").Value = $testData
$ws.Range("C14").Value = $testData
$ws.Range("D14").Value = $testData
$ws.Range("E14").Value = $testData
$ws.Rows.Item(14).RowHeight = 141.75

# --- Row 15: Post Conditions ---
$postCond = 'Message "Organization successfully registered"'
$ws.Range("B15").Value = $postCond
$ws.Range("C15").Value = $postCond
$ws.Range("D15").Value = $postCond
$ws.Range("E15").Value = $postCond
$ws.Rows.Item(15).RowHeight = 31.5

# --- Row 16: Expected Result ---
$expResult = "Successful register of an Organization"
$ws.Range("B16").Value = $expResult
$ws.Range("C16").Value = $expResult
$ws.Range("D16").Value = $expResult
$ws.Range("E16").Value = $expResult
$ws.Rows.Item(16).RowHeight = 31.5

# --- Row 17: Actual Result stays blank, but gets an explicit (empty) cell ---
$ws.Range("B17").Style = "Normal"

# --- Row 18: Status ---
$status = "Pass"
$ws.Range("B18").Value = $status
$ws.Range("C18").Value = $status
$ws.Range("D18").Value = $status
$ws.Range("E18").Value = $status

# --- Column widths (A:E -> 25.25 chars) ---
$ws.Columns.Item(1).ColumnWidth = 24.33
$ws.Columns.Item(2).ColumnWidth = 24.33
$ws.Columns.Item(3).ColumnWidth = 24.33
$ws.Columns.Item(4).ColumnWidth = 24.33
$ws.Columns.Item(5).ColumnWidth = 24.33

# --- View / selection state ---
# Sheet "Mover ficheiro" selection moves from E9 to B13
$ws2 = $wb.Worksheets.Item("Mover ficheiro")
$ws2.Range("B13").Select()

# "UC 1 - Registar Organização" becomes the active/selected tab, cursor at G11
$ws.Range("G11").Select()
$ws.Activate()
